$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 2518.9167  # H8
$ws.Cells.Item(8, 9).Value = 37.833332  # I8
$ws.Cells.Item(8, 11).Value = 113.499996  # K8
$ws.Cells.Item(8, 13).Value = 25.500004  # M8
$ws.Cells.Item(12, 8).Value = 534.25  # H12
$ws.Cells.Item(12, 9).Value = 534.25  # I12
$ws.Cells.Item(12, 11).Value = 534.25  # K12
$ws.Cells.Item(12, 13).Value = -364.25  # M12
$ws.Cells.Item(64, 8).Value = 7006.923  # H64
$ws.Cells.Item(64, 9).Value = 4846.8  # I64
$ws.Cells.Item(64, 11).Value = 4846.8  # K64
$ws.Cells.Item(64, 13).Value = -4598.8  # M64
$ws.Cells.Item(67, 8).Value = 7006.923  # H67
$ws.Cells.Item(67, 9).Value = 4846.8  # I67
$ws.Cells.Item(67, 11).Value = 4846.8  # K67
$ws.Cells.Item(67, 13).Value = -3988.8  # M67
$ws.Cells.Item(88, 8).Value = 1680.4286  # H88
$ws.Cells.Item(88, 10).Value = 1123.25  # J88
$ws.Cells.Item(88, 12).Value = 1123.25  # L88
$ws.Cells.Item(88, 14).Value = -1935.25  # N88
$ws.Cells.Item(91, 8).Value = 1680.4286  # H91
$ws.Cells.Item(91, 10).Value = 1123.25  # J91
$ws.Cells.Item(91, 12).Value = 1123.25  # L91
$ws.Cells.Item(91, 14).Value = -3931.25  # N91
$ws.Cells.Item(100, 8).Value = 3013.6667  # H100
$ws.Cells.Item(100, 9).Value = 3020.5  # I100
$ws.Cells.Item(100, 10).Value = 3000  # J100
$ws.Cells.Item(100, 11).Value = 3020.5  # K100
$ws.Cells.Item(100, 12).Value = 3000  # L100
$ws.Cells.Item(100, 13).Value = -2479.5  # M100
$ws.Cells.Item(100, 14).Value = -4082  # N100
$ws.Cells.Item(137, 8).Value = 6897743.5  # H137
$ws.Cells.Item(137, 9).Value = 8001168.5  # I137
$ws.Cells.Item(137, 11).Value = 24003505.5  # K137
$ws.Cells.Item(137, 13).Value = -24000955.5  # M137
$ws.Cells.Item(138, 8).Value = 5875.276  # H138
$ws.Cells.Item(138, 9).Value = 2898.3333  # I138
$ws.Cells.Item(138, 10).Value = 6651.8696  # J138
$ws.Cells.Item(138, 11).Value = 8694.999899999999  # K138
$ws.Cells.Item(138, 12).Value = 19955.6088  # L138
$ws.Cells.Item(138, 13).Value = -3554.999899999999  # M138
$ws.Cells.Item(138, 14).Value = -30235.6088  # N138
$ws.Cells.Item(141, 8).Value = 6552.0625  # H141
$ws.Cells.Item(141, 9).Value = 2672.8  # I141
$ws.Cells.Item(141, 11).Value = 8018.400000000001  # K141
$ws.Cells.Item(141, 13).Value = -2838.400000000001  # M141

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 470026.9  # H32
$ws.Cells.Item(32, 9).Value = 558173.4399999999  # I32
$ws.Cells.Item(32, 11).Value = 558173.4399999999  # K32
$ws.Cells.Item(32, 13).Value = -557886.4399999999  # M32
$ws.Cells.Item(45, 8).Value = 1872.75  # H45
$ws.Cells.Item(45, 9).Value = 1922.4  # I45
$ws.Cells.Item(45, 11).Value = 1922.4  # K45
$ws.Cells.Item(45, 13).Value = -1545.4  # M45
$ws.Cells.Item(61, 8).Value = 5978924.5  # H61
$ws.Cells.Item(61, 9).Value = 2459309.5  # I61
$ws.Cells.Item(61, 10).Value = 19001500  # J61
$ws.Cells.Item(61, 11).Value = 2459309.5  # K61
$ws.Cells.Item(61, 12).Value = 19001500  # L61
$ws.Cells.Item(61, 13).Value = -2459097.5  # M61
$ws.Cells.Item(61, 14).Value = -19001924  # N61
$ws.Cells.Item(130, 8).Value = 74996  # H130
$ws.Cells.Item(130, 10).Value = 74996  # J130
$ws.Cells.Item(130, 12).Value = 74996  # L130
$ws.Cells.Item(130, 14).Value = -85036  # N130
$ws.Cells.Item(132, 8).Value = 3227.3044  # H132
$ws.Cells.Item(132, 9).Value = 2011.4  # I132
$ws.Cells.Item(132, 10).Value = 11333.333  # J132
$ws.Cells.Item(132, 11).Value = 6034.200000000001  # K132
$ws.Cells.Item(132, 12).Value = 33999.999  # L132
$ws.Cells.Item(132, 13).Value = -3504.200000000001  # M132
$ws.Cells.Item(132, 14).Value = -39059.999  # N132
$ws.Cells.Item(136, 8).Value = 5978924.5  # H136
$ws.Cells.Item(136, 9).Value = 2459309.5  # I136
$ws.Cells.Item(136, 10).Value = 19001500  # J136
$ws.Cells.Item(136, 11).Value = 7377928.5  # K136
$ws.Cells.Item(136, 12).Value = 57004500  # L136
$ws.Cells.Item(136, 13).Value = -7375378.5  # M136
$ws.Cells.Item(136, 14).Value = -57009600  # N136
$ws.Cells.Item(140, 8).Value = 100132.14  # H140
$ws.Cells.Item(140, 10).Value = 50154.168  # J140
$ws.Cells.Item(140, 12).Value = 50154.168  # L140
$ws.Cells.Item(140, 14).Value = -60514.168  # N140

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2828.4285  # H86
$ws.Cells.Item(86, 9).Value = 1799.8334  # I86
$ws.Cells.Item(86, 11).Value = 1799.8334  # K86
$ws.Cells.Item(86, 13).Value = -676.8334  # M86
$ws.Cells.Item(89, 8).Value = 2828.4285  # H89
$ws.Cells.Item(89, 9).Value = 1799.8334  # I89
$ws.Cells.Item(89, 11).Value = 8999.166999999999  # K89
$ws.Cells.Item(89, 13).Value = -3383.166999999999  # M89
$ws.Cells.Item(105, 8).Value = 1630.125  # H105
$ws.Cells.Item(105, 9).Value = 1630.125  # I105
$ws.Cells.Item(105, 10).Value = 0  # J105
$ws.Cells.Item(105, 11).Value = 1630.125  # K105
$ws.Cells.Item(105, 12).Value = 0  # L105
$ws.Cells.Item(105, 13).Value = 116.875  # M105
$ws.Cells.Item(105, 14).ClearContents()  # N105
$ws.Cells.Item(134, 8).Value = 4833708  # H134
$ws.Cells.Item(134, 9).Value = 4067732.2  # I134
$ws.Cells.Item(134, 11).Value = 12203196.6  # K134
$ws.Cells.Item(134, 13).Value = -12200661.6  # M134

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1605.1666  # H22
$ws.Cells.Item(22, 9).Value = 968.73334  # I22
$ws.Cells.Item(22, 11).Value = 968.73334  # K22
$ws.Cells.Item(22, 13).Value = -618.73334  # M22
$ws.Cells.Item(105, 8).Value = 4677.0435  # H105
$ws.Cells.Item(105, 9).Value = 3819.85  # I105
$ws.Cells.Item(105, 10).Value = 10391.667  # J105
$ws.Cells.Item(105, 11).Value = 3819.85  # K105
$ws.Cells.Item(105, 12).Value = 10391.667  # L105
$ws.Cells.Item(105, 13).Value = -2072.85  # M105
$ws.Cells.Item(105, 14).Value = -13885.667  # N105
$ws.Cells.Item(134, 8).Value = 5839.222  # H134
$ws.Cells.Item(134, 9).Value = 4366.1  # I134
$ws.Cells.Item(134, 11).Value = 13098.3  # K134
$ws.Cells.Item(134, 13).Value = -10563.3  # M134

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 7644.8667  # H3
$ws.Cells.Item(3, 9).Value = 2852.5557  # I3
$ws.Cells.Item(3, 10).Value = 14833.333  # J3
$ws.Cells.Item(3, 11).Value = 8557.667099999999  # K3
$ws.Cells.Item(3, 12).Value = 44499.999  # L3
$ws.Cells.Item(3, 13).Value = -8445.667099999999  # M3
$ws.Cells.Item(3, 14).Value = -44723.999  # N3
$ws.Cells.Item(9, 8).Value = 8275.4  # H9
$ws.Cells.Item(9, 9).Value = 292  # I9
$ws.Cells.Item(9, 10).Value = 20250.5  # J9
$ws.Cells.Item(9, 11).Value = 876  # K9
$ws.Cells.Item(9, 12).Value = 60751.5  # L9
$ws.Cells.Item(9, 13).Value = -652  # M9
$ws.Cells.Item(9, 14).Value = -61199.5  # N9
$ws.Cells.Item(10, 8).Value = 2505  # H10
$ws.Cells.Item(10, 9).Value = 10  # I10
$ws.Cells.Item(10, 10).Value = 5000  # J10
$ws.Cells.Item(10, 11).Value = 30  # K10
$ws.Cells.Item(10, 12).Value = 15000  # L10
$ws.Cells.Item(10, 13).Value = 109  # M10
$ws.Cells.Item(10, 14).Value = -15278  # N10
$ws.Cells.Item(129, 8).Value = 627340.75  # H129
$ws.Cells.Item(129, 10).Value = 2459.4  # J129
$ws.Cells.Item(129, 12).Value = 7378.200000000001  # L129
$ws.Cells.Item(129, 14).Value = -17378.2  # N129
$ws.Cells.Item(137, 8).Value = 6204.9375  # H137
$ws.Cells.Item(137, 9).Value = 1150.4117  # I137
$ws.Cells.Item(137, 10).Value = 11933.4  # J137
$ws.Cells.Item(137, 11).Value = 3451.2351  # K137
$ws.Cells.Item(137, 12).Value = 35800.2  # L137
$ws.Cells.Item(137, 13).Value = 1648.7649  # M137
$ws.Cells.Item(137, 14).Value = -46000.2  # N137

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 5633  # H20
$ws.Cells.Item(20, 9).Value = 3782.8572  # I20
$ws.Cells.Item(20, 11).Value = 3782.8572  # K20
$ws.Cells.Item(20, 13).Value = -3537.8572  # M20
$ws.Cells.Item(49, 8).Value = 22079.166  # H49
$ws.Cells.Item(49, 9).Value = 21888  # I49
$ws.Cells.Item(49, 10).Value = 22117.4  # J49
$ws.Cells.Item(49, 11).Value = 21888  # K49
$ws.Cells.Item(49, 12).Value = 22117.4  # L49
$ws.Cells.Item(49, 13).Value = -21704  # M49
$ws.Cells.Item(49, 14).Value = -22485.4  # N49
$ws.Cells.Item(70, 8).Value = 42222.11  # H70
$ws.Cells.Item(70, 10).Value = 13333  # J70
$ws.Cells.Item(70, 12).Value = 13333  # L70
$ws.Cells.Item(70, 14).Value = -13873  # N70
$ws.Cells.Item(73, 8).Value = 42222.11  # H73
$ws.Cells.Item(73, 10).Value = 13333  # J73
$ws.Cells.Item(73, 12).Value = 13333  # L73
$ws.Cells.Item(73, 14).Value = -15205  # N73
$ws.Cells.Item(134, 8).Value = 0  # H134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 12).Value = 0  # L134
$ws.Cells.Item(134, 14).ClearContents()  # N134

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 2050002  # H2
$ws.Cells.Item(2, 9).Value = 2050002  # I2
$ws.Cells.Item(2, 11).Value = 2050002  # K2
$ws.Cells.Item(2, 13).Value = -2049890  # M2
$ws.Cells.Item(22, 8).Value = 3124.35  # H22
$ws.Cells.Item(22, 9).Value = 2441.8  # I22
$ws.Cells.Item(22, 10).Value = 3806.9  # J22
$ws.Cells.Item(22, 11).Value = 2441.8  # K22
$ws.Cells.Item(22, 12).Value = 3806.9  # L22
$ws.Cells.Item(22, 13).Value = -2146.8  # M22
$ws.Cells.Item(22, 14).Value = -4396.9  # N22
$ws.Cells.Item(27, 8).Value = 3124.35  # H27
$ws.Cells.Item(27, 9).Value = 2441.8  # I27
$ws.Cells.Item(27, 10).Value = 3806.9  # J27
$ws.Cells.Item(27, 11).Value = 2441.8  # K27
$ws.Cells.Item(27, 12).Value = 3806.9  # L27
$ws.Cells.Item(27, 13).Value = -2334.8  # M27
$ws.Cells.Item(27, 14).Value = -4020.9  # N27
$ws.Cells.Item(46, 8).Value = 3149.375  # H46
$ws.Cells.Item(46, 9).Value = 1800  # I46
$ws.Cells.Item(46, 10).Value = 3239.3333  # J46
$ws.Cells.Item(46, 11).Value = 1800  # K46
$ws.Cells.Item(46, 12).Value = 3239.3333  # L46
$ws.Cells.Item(46, 13).Value = -1612  # M46
$ws.Cells.Item(46, 14).Value = -3615.3333  # N46
$ws.Cells.Item(55, 8).Value = 2336.9  # H55
$ws.Cells.Item(55, 9).Value = 1174.8  # I55
$ws.Cells.Item(55, 11).Value = 1174.8  # K55
$ws.Cells.Item(55, 13).Value = -1001.8  # M55
$ws.Cells.Item(68, 8).Value = 6636.4  # H68
$ws.Cells.Item(68, 10).Value = 9499.5  # J68
$ws.Cells.Item(68, 12).Value = 9499.5  # L68
$ws.Cells.Item(68, 14).Value = -10997.5  # N68
$ws.Cells.Item(71, 8).Value = 6636.4  # H71
$ws.Cells.Item(71, 10).Value = 9499.5  # J71
$ws.Cells.Item(71, 12).Value = 47497.5  # L71
$ws.Cells.Item(71, 14).Value = -54985.5  # N71
$ws.Cells.Item(93, 8).Value = 3653.889  # H93
$ws.Cells.Item(93, 9).Value = 1518.4  # I93
$ws.Cells.Item(93, 11).Value = 1518.4  # K93
$ws.Cells.Item(93, 13).Value = -270.4000000000001  # M93
$ws.Cells.Item(133, 8).Value = 85763.664  # H133
$ws.Cells.Item(133, 10).Value = 85763.664  # J133
$ws.Cells.Item(133, 12).Value = 85763.664  # L133
$ws.Cells.Item(133, 14).Value = -90823.664  # N133

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3982.2778  # H81
$ws.Cells.Item(81, 9).Value = 3673.875  # I81
$ws.Cells.Item(81, 10).Value = 4229  # J81
$ws.Cells.Item(81, 11).Value = 7347.75  # K81
$ws.Cells.Item(81, 12).Value = 8458  # L81
$ws.Cells.Item(81, 13).Value = -6286.75  # M81
$ws.Cells.Item(81, 14).Value = -10580  # N81
$ws.Cells.Item(84, 8).Value = 3982.2778  # H84
$ws.Cells.Item(84, 9).Value = 3673.875  # I84
$ws.Cells.Item(84, 10).Value = 4229  # J84
$ws.Cells.Item(84, 11).Value = 36738.75  # K84
$ws.Cells.Item(84, 12).Value = 42290  # L84
$ws.Cells.Item(84, 13).Value = -31434.75  # M84
$ws.Cells.Item(84, 14).Value = -52898  # N84
$ws.Cells.Item(96, 8).Value = 1729.3334  # H96
$ws.Cells.Item(96, 9).Value = 1729.3334  # I96
$ws.Cells.Item(96, 11).Value = 1729.3334  # K96
$ws.Cells.Item(96, 13).Value = -356.3334  # M96
$ws.Cells.Item(132, 8).Value = 3625643.2  # H132
$ws.Cells.Item(132, 9).Value = 4275987  # I132
$ws.Cells.Item(132, 11).Value = 12827961  # K132
$ws.Cells.Item(132, 13).Value = -12825431  # M132
